# SO_Test.xlsx edit: apply number formatting to DEPOSIT/PRICE columns and
# append four new sales-order rows (8-11) below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Number-format the existing DEPOSIT (G) / PRICE (H) values -------------
# Whole numbers -> "#,##0"; values with a fractional part -> "#,##0.00"
$moneyCells = @(
    @{ Cell = "G2"; Fmt = "#,##0" },
    @{ Cell = "H2"; Fmt = "#,##0.00" },
    @{ Cell = "G3"; Fmt = "#,##0" },
    @{ Cell = "H3"; Fmt = "#,##0" },
    @{ Cell = "G5"; Fmt = "#,##0" },
    @{ Cell = "H5"; Fmt = "#,##0.00" },
    @{ Cell = "G7"; Fmt = "#,##0" },
    @{ Cell = "H7"; Fmt = "#,##0" }
)
foreach ($mc in $moneyCells) {
    $rng = $ws.Range($mc.Cell)
    $rng.Font.Name = "Calibri"
    $rng.NumberFormat = $mc.Fmt
    $rng.HorizontalAlignment = -4152
}

# --- Row heights -------------------------------------------------------
foreach ($r in 2..7) {
    $ws.Rows.Item($r).RowHeight = 18.75
}

# --- Pad rows 4 and 6 (the blank separator rows) with formatted empty cells
foreach ($r in @(4, 6)) {
    $ws.Range("A$r`:F$r").Font.Name = "Calibri"
    $ws.Range("I$r`:K$r").Font.Name = "Calibri"
    $g = $ws.Range("G$r")
    $g.NumberFormat = "#,##0"
    $g.HorizontalAlignment = -4152
    $h = $ws.Range("H$r")
    $h.NumberFormat = "#,##0"
    $h.HorizontalAlignment = -4152
}

# --- Append new rows 8-11 -----------------------------------------------
$newRows = @(
    @("01/12/2024","SO240112002","Ashley","9157994875","ATTILA","ABOUT THAT LIFE",25,65,"AMS","CD","Ashley"),
    @("01/12/2024","SO240112003","Ashley","9157994875","Chicago","EH",15,74,"AMA","DVD","Ashley"),
    @("01/12/2024","SO240112004","A","9157994875","666","999",12,999,"AMS","LP","A"),
    @("01/12/2024","SO240112005","ASHLEY","9157994875","ASAASASASAS","OSDFGJKH",25,3333,"AMS","BLU-RAY","A")
)

$rowNum = 8
foreach ($rowVals in $newRows) {
    $colLetters = @("A","B","C","D","E","F","G","H","I","J","K")
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $addr = "$($colLetters[$i])$rowNum"
        $val = $rowVals[$i]
        $cell = $ws.Range($addr)
        if ($i -eq 6 -or $i -eq 7) {
            # DEPOSIT / PRICE -> numeric
            $cell.Value = $val
        } else {
            # Everything else stored as literal text (dates, phone numbers and
            # numeric-looking artist/title strings must NOT become real numbers)
            $cell.NumberFormat = "@"
            $cell.Value = [string]$val
        }
    }
    $rowNum++
}
